# Auto-generated PowerShell COM-interop script
# Applies shape position/size (and two flip) changes per the target diff
# for docs/diagrams/SortSequenceDiagram.pptx (slide 1).

function Get-ShapeById {
    param($slide, $targetId)
    $count = $slide.Shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$updates = @(
    @{ Id = 47; Left = 726.7080078125; Top = -11.112283706665039; Width = 176.0245819091797; Height = 593.11181640625 }  # Rectangle 65
    @{ Id = 262; Left = 848.7039794921875; Top = 488.6813659667969 }  # Straight Connector 261
    @{ Id = 263; Left = 842.8659057617188; Top = 504.0000915527344; Width = 11.676141738891602; Height = 18.0 }  # Rectangle 262
    @{ Id = 48; Left = -62.66653823852539; Top = -11.1123628616333; Width = 784.6602783203125; Height = 593.1119384765625 }  # Rectangle 65
    @{ Id = 57; Width = 0.0; Height = 565.5789184570312 }  # Straight Connector 56
    @{ Id = 58; Width = 13.425827026367188; Height = 512.7216796875 }  # Rectangle 57
    @{ Id = 87; Left = -68.39409637451172; Top = 558.0432739257812; Width = 117.12669372558594; Height = 0.6680315136909485 }  # Straight Arrow Connector 86
    @{ Id = 89; Left = 603.9393310546875; Top = 329.2424621582031 }  # TextBox 88
    @{ Id = 94; Left = 731.7294921875; Top = 301.6834716796875 }  # Rectangle 62
    @{ Id = 97; Left = 591.4352416992188; Top = 344.7188415527344 }  # Straight Arrow Connector 96
    @{ Id = 98; Left = 662.255859375; Top = 387.5263977050781 }  # TextBox 97
    @{ Id = 100; Left = 54.410160064697266; Top = 290.12969970703125; Width = 525.5677490234375; Height = 2.044015884399414 }  # Straight Arrow Connector 99
    @{ Id = 102; Left = 54.410160064697266; Top = 540.0; Width = 530.3467407226562; Height = 2.044015884399414; VerticalFlip = $true }  # Straight Arrow Connector 101
    @{ Id = 108; Left = 764.5789184570312; Top = 522.0; Width = 79.2957534790039; Height = 0.29937008023262024 }  # Straight Arrow Connector 107
    @{ Id = 114; Left = 762.9971313476562; Top = 505.5313415527344; Width = 79.2957534790039; Height = 0.29937008023262024 }  # Straight Arrow Connector 113
    @{ Id = 116; Left = 263.99127197265625; Top = 526.91650390625 }  # TextBox 115
    @{ Id = 36; Left = 548.9144897460938; Top = 365.2424621582031 }  # Rectangle 35
    @{ Id = 38; Left = 549.0804443359375; Top = 365.3806457519531 }  # Snip Single Corner Rectangle 37
    @{ Id = 39; Left = 548.1920776367188; Top = 362.96630859375 }  # TextBox 38
    @{ Id = 124; Left = 612.8148193359375; Top = 363.55096435546875 }  # TextBox 123
    @{ Id = 125; Left = 585.22900390625; Top = 227.0550537109375; Width = 0.984094500541687; Height = 334.2403259277344 }  # Straight Connector 124
    @{ Id = 45; Left = 548.9144287109375; Top = 443.2424621582031 }  # Straight Connector 44
    @{ Id = 127; Left = 579.3953247070312; Top = 289.1993103027344; Width = 10.723150253295898; Height = 250.80072021484375 }  # Rectangle 126
    @{ Id = 128; Left = 549.9385986328125; Top = 198.0 }  # Rectangle 62
    @{ Id = 129; Left = 640.8175048828125; Top = 441.85498046875 }  # TextBox 128
    @{ Id = 143; Left = 585.4998168945312; Top = 456.718994140625 }  # TextBox 142
    @{ Id = 179; Left = 411.6812744140625; Top = 205.9118194580078 }  # TextBox 178
    @{ Id = 180; Left = 384.3046569824219; Top = 224.52394104003906; Width = 164.88174438476562; Height = 0.0 }  # Straight Arrow Connector 179
    @{ Id = 181; Left = 384.3046569824219; Top = 247.51678466796875; Width = 201.67323303222656; Height = 0.0 }  # Straight Arrow Connector 180
    @{ Id = 189; Left = 579.9779052734375; Top = 231.67103576660156 }  # Rectangle 188
    @{ Id = 200; Left = -42.50047302246094; Top = 540.3053588867188 }  # TextBox 199
    @{ Id = 210; Left = 648.8365478515625; Top = 282.7679748535156 }  # Rectangle 62
    @{ Id = 211; Left = 676.7724609375; Top = 308.1542663574219 }  # Rectangle 210
    @{ Id = 212; Left = 590.1397094726562; Top = 300.62811279296875 }  # Straight Arrow Connector 211
    @{ Id = 214; Left = 591.9779052734375; Top = 324.0 }  # Straight Arrow Connector 213
    @{ Id = 216; Left = 587.9201049804688; Top = 273.2889099121094 }  # TextBox 215
    @{ Id = 217; Left = 609.1661987304688; Top = 310.6903381347656 }  # TextBox 216
    @{ Id = 223; Left = 611.069091796875; Top = 380.75433349609375 }  # TextBox 222
    @{ Id = 224; Left = 588.193115234375; Top = 386.5839538574219 }  # Elbow Connector 223
    @{ Id = 225; Left = 583.2689208984375; Top = 430.6715087890625 }  # Elbow Connector 224
    @{ Id = 226; Left = 585.7506103515625; Top = 392.18402099609375 }  # Rectangle 225
    @{ Id = 233; Left = 591.9779052734375; Top = 467.3429260253906 }  # Straight Arrow Connector 232
    @{ Id = 95; Left = 757.9561767578125; Top = 325.5315856933594; Width = 0.0; Height = 241.4189910888672; HorizontalFlip = $false }  # Straight Connector 94
    @{ Id = 96; Left = 751.3209838867188; Top = 338.6395568847656; Width = 11.676141738891602; Height = 192.7523651123047 }  # Rectangle 95
    @{ Id = 103; Left = 591.4352416992188; Top = 355.5622253417969 }  # Straight Arrow Connector 102
    @{ Id = 237; Left = 589.81201171875; Top = 402.1866149902344 }  # TextBox 236
    @{ Id = 238; Left = 595.0444946289062; Top = 413.1108703613281 }  # Straight Arrow Connector 237
    @{ Id = 252; Left = 595.0444946289062; Top = 420.9042663574219 }  # Straight Arrow Connector 251
    @{ Id = 253; Left = 591.4971313476562; Top = 476.99237060546875 }  # Straight Arrow Connector 252
    @{ Id = 254; Left = 584.9360961914062; Top = 488.6813659667969 }  # TextBox 253
    @{ Id = 255; Left = 591.4141845703125; Top = 499.3052978515625 }  # Straight Arrow Connector 254
    @{ Id = 256; Left = 590.92724609375; Top = 531.2661743164062 }  # Straight Arrow Connector 255
    @{ Id = 261; Left = 807.2376098632812; Top = 456.1073303222656 }  # Rectangle 62
    @{ Id = 264; Left = 764.9174194335938; Top = 492.10552978515625 }  # TextBox 263
    @{ Id = 284; Left = 573.8226318359375; Top = 546.0 }  # TextBox 283
)

foreach ($u in $updates) {
    $sh = Get-ShapeById $s $u.Id
    if ($null -eq $sh) {
        Write-Host "WARNING: shape id" $u.Id "not found"
        continue
    }
    if ($u.ContainsKey('Width')) { $sh.Width = $u.Width }
    if ($u.ContainsKey('Height')) { $sh.Height = $u.Height }
    if ($u.ContainsKey('Left')) { $sh.Left = $u.Left }
    if ($u.ContainsKey('Top')) { $sh.Top = $u.Top }
    if ($u.ContainsKey('VerticalFlip')) { $sh.VerticalFlip = $u.VerticalFlip }
    if ($u.ContainsKey('HorizontalFlip')) { $sh.HorizontalFlip = $u.HorizontalFlip }
}

Write-Host "Applied" $updates.Count "shape updates"
